# Auto-generated Excel COM-interop script
# Applies numeric updates to columns H-N across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 2583.5  # H6: 2966.4285 -> 2583.5
$ws.Cells.Item(6, 9).Value = 4033.6  # I6: 4132.6 -> 4033.6
$ws.Cells.Item(6, 10).Value = 166.66667  # J6: 51 -> 166.66667
$ws.Cells.Item(6, 11).Value = 12100.8  # K6: 12397.8 -> 12100.8
$ws.Cells.Item(6, 12).Value = 500.00001  # L6: 153 -> 500.00001
$ws.Cells.Item(6, 13).Value = -11988.8  # M6: -12285.8 -> -11988.8
$ws.Cells.Item(6, 14).Value = -724.00001  # N6: -377 -> -724.00001
$ws.Cells.Item(11, 8).Value = 585.63635  # H11: 577.6667 -> 585.63635
$ws.Cells.Item(11, 9).Value = 585.63635  # I11: 577.6667 -> 585.63635
$ws.Cells.Item(11, 11).Value = 585.63635  # K11: 577.6667 -> 585.63635
$ws.Cells.Item(11, 13).Value = -445.63635  # M11: -437.6667 -> -445.63635
$ws.Cells.Item(32, 8).Value = 954.8889  # H32: 924.4167 -> 954.8889
$ws.Cells.Item(32, 9).Value = 800  # I32: 824.75 -> 800
$ws.Cells.Item(32, 10).Value = 999.1429000000001  # J32: 974.25 -> 999.1429000000001
$ws.Cells.Item(32, 11).Value = 800  # K32: 824.75 -> 800
$ws.Cells.Item(32, 12).Value = 999.1429000000001  # L32: 974.25 -> 999.1429000000001
$ws.Cells.Item(32, 13).Value = -474  # M32: -498.75 -> -474
$ws.Cells.Item(32, 14).Value = -1651.1429  # N32: -1626.25 -> -1651.1429
$ws.Cells.Item(40, 8).Value = 1490  # H40: 1495 -> 1490
$ws.Cells.Item(40, 9).Value = 1490  # I40: 0 -> 1490
$ws.Cells.Item(40, 10).Value = 0  # J40: 1495 -> 0
$ws.Cells.Item(40, 11).Value = 1490  # K40: 0 -> 1490
$ws.Cells.Item(40, 12).Value = 0  # L40: 1495 -> 0
$ws.Cells.Item(40, 13).Value = -1315  # M40: None -> -1315
$ws.Cells.Item(40, 14).Value = $null  # N40: -1845 -> (removed)
$ws.Cells.Item(43, 8).Value = 2395.6365  # H43: 56168 -> 2395.6365
$ws.Cells.Item(43, 10).Value = 2870.4  # J43: 96794 -> 2870.4
$ws.Cells.Item(43, 12).Value = 2870.4  # L43: 96794 -> 2870.4
$ws.Cells.Item(43, 14).Value = -3008.4  # N43: -96932 -> -3008.4
$ws.Cells.Item(58, 8).Value = 935.4286  # H58: 1642.375 -> 935.4286
$ws.Cells.Item(58, 9).Value = 591.3333  # I58: 629.8 -> 591.3333
$ws.Cells.Item(58, 10).Value = 3000  # J58: 3330 -> 3000
$ws.Cells.Item(58, 11).Value = 1773.9999  # K58: 1889.4 -> 1773.9999
$ws.Cells.Item(58, 12).Value = 9000  # L58: 9990 -> 9000
$ws.Cells.Item(58, 13).Value = -1623.9999  # M58: -1739.4 -> -1623.9999
$ws.Cells.Item(58, 14).Value = -9300  # N58: -10290 -> -9300
$ws.Cells.Item(106, 8).Value = 8599  # H106: 3591.3333 -> 8599
$ws.Cells.Item(106, 9).Value = 8599  # I106: 3591.3333 -> 8599
$ws.Cells.Item(106, 11).Value = 8599  # K106: 3591.3333 -> 8599
$ws.Cells.Item(106, 13).Value = -7968  # M106: -2960.3333 -> -7968
$ws.Cells.Item(116, 8).Value = 5447.1665  # H116: 6030.5 -> 5447.1665
$ws.Cells.Item(116, 9).Value = 2000  # I116: 5500 -> 2000
$ws.Cells.Item(116, 11).Value = 2000  # K116: 5500 -> 2000
$ws.Cells.Item(116, 13).Value = 1442  # M116: -2058 -> 1442

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1016.82355  # H2: 1052.4375 -> 1016.82355
$ws.Cells.Item(2, 9).Value = 1036.625  # I2: 1075.9333 -> 1036.625
$ws.Cells.Item(2, 11).Value = 1036.625  # K2: 1075.9333 -> 1036.625
$ws.Cells.Item(2, 13).Value = -923.625  # M2: -962.9332999999999 -> -923.625
$ws.Cells.Item(3, 8).Value = 7000.6  # H3: 9629 -> 7000.6
$ws.Cells.Item(3, 9).Value = 5750.75  # I3: 9629 -> 5750.75
$ws.Cells.Item(3, 10).Value = 12000  # J3: 0 -> 12000
$ws.Cells.Item(3, 11).Value = 5750.75  # K3: 9629 -> 5750.75
$ws.Cells.Item(3, 12).Value = 12000  # L3: 0 -> 12000
$ws.Cells.Item(3, 13).Value = -5635.75  # M3: -9514 -> -5635.75
$ws.Cells.Item(3, 14).Value = -12230  # N3: None -> -12230
$ws.Cells.Item(45, 8).Value = 1864.1818  # H45: 1901.1 -> 1864.1818
$ws.Cells.Item(45, 9).Value = 1800.6666  # I45: 1838.875 -> 1800.6666
$ws.Cells.Item(45, 11).Value = 1800.6666  # K45: 1838.875 -> 1800.6666
$ws.Cells.Item(45, 13).Value = -1423.6666  # M45: -1461.875 -> -1423.6666
$ws.Cells.Item(74, 8).Value = 1263.4  # H74: 1337.2222 -> 1263.4
$ws.Cells.Item(74, 9).Value = 1263.4  # I74: 1337.2222 -> 1263.4
$ws.Cells.Item(74, 11).Value = 1263.4  # K74: 1337.2222 -> 1263.4
$ws.Cells.Item(74, 13).Value = -389.4000000000001  # M74: -463.2221999999999 -> -389.4000000000001
$ws.Cells.Item(77, 8).Value = 1263.4  # H77: 1337.2222 -> 1263.4
$ws.Cells.Item(77, 9).Value = 1263.4  # I77: 1337.2222 -> 1263.4
$ws.Cells.Item(77, 11).Value = 6317  # K77: 6686.111 -> 6317
$ws.Cells.Item(77, 13).Value = -1949  # M77: -2318.111 -> -1949
$ws.Cells.Item(116, 8).Value = 1016.82355  # H116: 1052.4375 -> 1016.82355
$ws.Cells.Item(116, 9).Value = 1036.625  # I116: 1075.9333 -> 1036.625
$ws.Cells.Item(116, 11).Value = 1036.625  # K116: 1075.9333 -> 1036.625
$ws.Cells.Item(116, 13).Value = 1257.375  # M116: 1218.0667 -> 1257.375
$ws.Cells.Item(122, 8).Value = 2378.3635  # H122: 2242.5833 -> 2378.3635
$ws.Cells.Item(122, 9).Value = 2378.3635  # I122: 2242.5833 -> 2378.3635
$ws.Cells.Item(122, 11).Value = 7135.0905  # K122: 6727.749899999999 -> 7135.0905
$ws.Cells.Item(122, 13).Value = -4685.0905  # M122: -4277.749899999999 -> -4685.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1016.82355  # H3: 1052.4375 -> 1016.82355
$ws.Cells.Item(3, 9).Value = 1036.625  # I3: 1075.9333 -> 1036.625
$ws.Cells.Item(3, 11).Value = 1036.625  # K3: 1075.9333 -> 1036.625
$ws.Cells.Item(3, 13).Value = -922.625  # M3: -961.9332999999999 -> -922.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2000  # H62: 0 -> 2000
$ws.Cells.Item(62, 10).Value = 2000  # J62: 0 -> 2000
$ws.Cells.Item(62, 12).Value = 2000  # L62: 0 -> 2000
$ws.Cells.Item(62, 14).Value = -3248  # N62: None -> -3248
$ws.Cells.Item(65, 8).Value = 2000  # H65: 0 -> 2000
$ws.Cells.Item(65, 10).Value = 2000  # J65: 0 -> 2000
$ws.Cells.Item(65, 12).Value = 10000  # L65: 0 -> 10000
$ws.Cells.Item(65, 14).Value = -16240  # N65: None -> -16240
$ws.Cells.Item(99, 8).Value = 8374.5  # H99: 8749 -> 8374.5
$ws.Cells.Item(99, 10).Value = 8000  # J99: 0 -> 8000
$ws.Cells.Item(99, 12).Value = 8000  # L99: 0 -> 8000
$ws.Cells.Item(99, 14).Value = -10996  # N99: None -> -10996
$ws.Cells.Item(107, 8).Value = 1165.75  # H107: 1419.6 -> 1165.75
$ws.Cells.Item(107, 9).Value = 1226.2858  # I107: 1419.6 -> 1226.2858
$ws.Cells.Item(107, 10).Value = 742  # J107: 0 -> 742
$ws.Cells.Item(107, 11).Value = 1226.2858  # K107: 1419.6 -> 1226.2858
$ws.Cells.Item(107, 12).Value = 742  # L107: 0 -> 742
$ws.Cells.Item(107, 13).Value = 693.7141999999999  # M107: 500.4000000000001 -> 693.7141999999999
$ws.Cells.Item(107, 14).Value = -4582  # N107: None -> -4582
$ws.Cells.Item(122, 8).Value = 1716.6666  # H122: 2600 -> 1716.6666
$ws.Cells.Item(122, 9).Value = 1716.6666  # I122: 2600 -> 1716.6666
$ws.Cells.Item(122, 11).Value = 5149.9998  # K122: 7800 -> 5149.9998
$ws.Cells.Item(122, 13).Value = -2699.9998  # M122: -5350 -> -2699.9998
$ws.Cells.Item(126, 8).Value = 8374.5  # H126: 8749 -> 8374.5
$ws.Cells.Item(126, 10).Value = 8000  # J126: 0 -> 8000
$ws.Cells.Item(126, 12).Value = 24000  # L126: 0 -> 24000
$ws.Cells.Item(126, 14).Value = -28940  # N126: None -> -28940

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(16, 8).Value = 1000  # H16: 640.2857 -> 1000
$ws.Cells.Item(16, 9).Value = 0  # I16: 41.5 -> 0
$ws.Cells.Item(16, 10).Value = 1000  # J16: 879.8 -> 1000
$ws.Cells.Item(16, 11).Value = 0  # K16: 124.5 -> 0
$ws.Cells.Item(16, 12).Value = 3000  # L16: 2639.4 -> 3000
$ws.Cells.Item(16, 13).Value = $null  # M16: 48.5 -> (removed)
$ws.Cells.Item(16, 14).Value = -3346  # N16: -2985.4 -> -3346
$ws.Cells.Item(22, 8).Value = 2591.875  # H22: 2542.353 -> 2591.875
$ws.Cells.Item(22, 10).Value = 2747.7693  # J22: 2676.5 -> 2747.7693
$ws.Cells.Item(22, 12).Value = 8243.3079  # L22: 8029.5 -> 8243.3079
$ws.Cells.Item(22, 14).Value = -8581.3079  # N22: -8367.5 -> -8581.3079
$ws.Cells.Item(27, 8).Value = 2591.875  # H27: 2542.353 -> 2591.875
$ws.Cells.Item(27, 10).Value = 2747.7693  # J27: 2676.5 -> 2747.7693
$ws.Cells.Item(27, 12).Value = 8243.3079  # L27: 8029.5 -> 8243.3079
$ws.Cells.Item(27, 14).Value = -8447.3079  # N27: -8233.5 -> -8447.3079
$ws.Cells.Item(63, 8).Value = 0  # H63: 2357.5 -> 0
$ws.Cells.Item(63, 9).Value = 0  # I63: 2222 -> 0
$ws.Cells.Item(63, 10).Value = 0  # J63: 2493 -> 0
$ws.Cells.Item(63, 11).Value = 0  # K63: 6666 -> 0
$ws.Cells.Item(63, 12).Value = 0  # L63: 7479 -> 0
$ws.Cells.Item(63, 13).Value = $null  # M63: -5917 -> (removed)
$ws.Cells.Item(63, 14).Value = $null  # N63: -8977 -> (removed)
$ws.Cells.Item(66, 8).Value = 0  # H66: 2357.5 -> 0
$ws.Cells.Item(66, 9).Value = 0  # I66: 2222 -> 0
$ws.Cells.Item(66, 10).Value = 0  # J66: 2493 -> 0
$ws.Cells.Item(66, 11).Value = 0  # K66: 19998 -> 0
$ws.Cells.Item(66, 12).Value = 0  # L66: 22437 -> 0
$ws.Cells.Item(66, 13).Value = $null  # M66: -16254 -> (removed)
$ws.Cells.Item(66, 14).Value = $null  # N66: -29925 -> (removed)
$ws.Cells.Item(75, 8).Value = 2222  # H75: 0 -> 2222
$ws.Cells.Item(75, 10).Value = 2222  # J75: 0 -> 2222
$ws.Cells.Item(75, 12).Value = 6666  # L75: 0 -> 6666
$ws.Cells.Item(75, 14).Value = -8662  # N75: None -> -8662
$ws.Cells.Item(78, 8).Value = 2222  # H78: 0 -> 2222
$ws.Cells.Item(78, 10).Value = 2222  # J78: 0 -> 2222
$ws.Cells.Item(78, 12).Value = 19998  # L78: 0 -> 19998
$ws.Cells.Item(78, 14).Value = -29982  # N78: None -> -29982
$ws.Cells.Item(118, 8).Value = 1495  # H118: 1500 -> 1495
$ws.Cells.Item(118, 9).Value = 1495  # I118: 1500 -> 1495
$ws.Cells.Item(118, 11).Value = 4485  # K118: 4500 -> 4485
$ws.Cells.Item(118, 13).Value = -3242  # M118: -3257 -> -3242
$ws.Cells.Item(119, 8).Value = 0  # H119: 129 -> 0
$ws.Cells.Item(119, 9).Value = 0  # I119: 129 -> 0
$ws.Cells.Item(119, 11).Value = 0  # K119: 387 -> 0
$ws.Cells.Item(119, 13).Value = $null  # M119: 4451 -> (removed)
$ws.Cells.Item(129, 8).Value = 300  # H129: 359 -> 300
$ws.Cells.Item(129, 10).Value = 500  # J129: 488.5 -> 500
$ws.Cells.Item(129, 12).Value = 1500  # L129: 1465.5 -> 1500
$ws.Cells.Item(129, 14).Value = -11500  # N129: -11465.5 -> -11500
$ws.Cells.Item(131, 8).Value = 1399.6666  # H131: 1199.6 -> 1399.6666
$ws.Cells.Item(131, 9).Value = 999  # I131: 932.6667 -> 999
$ws.Cells.Item(131, 11).Value = 2997  # K131: 2798.0001 -> 2997
$ws.Cells.Item(131, 13).Value = 2043  # M131: 2241.9999 -> 2043
$ws.Cells.Item(138, 8).Value = 6888.4443  # H138: 6749 -> 6888.4443
$ws.Cells.Item(138, 10).Value = 7466.8  # J138: 8400.4 -> 7466.8
$ws.Cells.Item(138, 12).Value = 22400.4  # L138: 25201.2 -> 22400.4
$ws.Cells.Item(138, 14).Value = -32680.4  # N138: -35481.2 -> -32680.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 2000  # H41: 2391.6667 -> 2000
$ws.Cells.Item(41, 9).Value = 2000  # I41: 4350 -> 2000
$ws.Cells.Item(41, 10).Value = 0  # J41: 2000 -> 0
$ws.Cells.Item(41, 11).Value = 2000  # K41: 4350 -> 2000
$ws.Cells.Item(41, 12).Value = 0  # L41: 2000 -> 0
$ws.Cells.Item(41, 13).Value = -1645  # M41: -3995 -> -1645
$ws.Cells.Item(41, 14).Value = $null  # N41: -2710 -> (removed)
$ws.Cells.Item(52, 8).Value = 30000  # H52: 0 -> 30000
$ws.Cells.Item(52, 10).Value = 30000  # J52: 0 -> 30000
$ws.Cells.Item(52, 12).Value = 30000  # L52: 0 -> 30000
$ws.Cells.Item(52, 14).Value = -30518  # N52: None -> -30518
$ws.Cells.Item(70, 8).Value = 7498.25  # H70: 8331.333000000001 -> 7498.25
$ws.Cells.Item(70, 9).Value = 4999.5  # I70: 5000 -> 4999.5
$ws.Cells.Item(70, 11).Value = 4999.5  # K70: 5000 -> 4999.5
$ws.Cells.Item(70, 13).Value = -4729.5  # M70: -4730 -> -4729.5
$ws.Cells.Item(73, 8).Value = 7498.25  # H73: 8331.333000000001 -> 7498.25
$ws.Cells.Item(73, 9).Value = 4999.5  # I73: 5000 -> 4999.5
$ws.Cells.Item(73, 11).Value = 4999.5  # K73: 5000 -> 4999.5
$ws.Cells.Item(73, 13).Value = -4063.5  # M73: -4064 -> -4063.5
$ws.Cells.Item(96, 8).Value = 29930  # H96: 0 -> 29930
$ws.Cells.Item(96, 10).Value = 29930  # J96: 0 -> 29930
$ws.Cells.Item(96, 12).Value = 29930  # L96: 0 -> 29930
$ws.Cells.Item(96, 14).Value = -35422  # N96: None -> -35422
$ws.Cells.Item(102, 8).Value = 1966.6666  # H102: 824.3333 -> 1966.6666
$ws.Cells.Item(102, 9).Value = 1450  # I102: 524.0909 -> 1450
$ws.Cells.Item(102, 10).Value = 3000  # J102: 1650 -> 3000
$ws.Cells.Item(102, 11).Value = 1450  # K102: 524.0909 -> 1450
$ws.Cells.Item(102, 12).Value = 3000  # L102: 1650 -> 3000
$ws.Cells.Item(102, 13).Value = 172  # M102: 1097.9091 -> 172
$ws.Cells.Item(102, 14).Value = -6244  # N102: -4894 -> -6244
$ws.Cells.Item(122, 8).Value = 6253050.5  # H122: 6582063.5 -> 6253050.5
$ws.Cells.Item(122, 10).Value = 2758.4  # J122: 2998.25 -> 2758.4
$ws.Cells.Item(122, 12).Value = 8275.200000000001  # L122: 8994.75 -> 8275.200000000001
$ws.Cells.Item(122, 14).Value = -13175.2  # N122: -13894.75 -> -13175.2
$ws.Cells.Item(132, 8).Value = 1671.75  # H132: 0 -> 1671.75
$ws.Cells.Item(132, 9).Value = 1671.75  # I132: 0 -> 1671.75
$ws.Cells.Item(132, 11).Value = 5015.25  # K132: 0 -> 5015.25
$ws.Cells.Item(132, 13).Value = -2485.25  # M132: None -> -2485.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2005.8667  # H7: 2007.3077 -> 2005.8667
$ws.Cells.Item(7, 9).Value = 1807.6364  # I7: 2089.1 -> 1807.6364
$ws.Cells.Item(7, 10).Value = 2551  # J7: 1734.6666 -> 2551
$ws.Cells.Item(7, 11).Value = 1807.6364  # K7: 2089.1 -> 1807.6364
$ws.Cells.Item(7, 12).Value = 2551  # L7: 1734.6666 -> 2551
$ws.Cells.Item(7, 13).Value = -1695.6364  # M7: -1977.1 -> -1695.6364
$ws.Cells.Item(7, 14).Value = -2775  # N7: -1958.6666 -> -2775
$ws.Cells.Item(46, 8).Value = 3166.6667  # H46: 3857.1428 -> 3166.6667
$ws.Cells.Item(46, 9).Value = 2555.5557  # I46: 3000 -> 2555.5557
$ws.Cells.Item(46, 11).Value = 2555.5557  # K46: 3000 -> 2555.5557
$ws.Cells.Item(46, 13).Value = -2367.5557  # M46: -2812 -> -2367.5557
$ws.Cells.Item(55, 8).Value = 374.5  # H55: 274.5 -> 374.5
$ws.Cells.Item(55, 10).Value = 600  # J55: 200 -> 600
$ws.Cells.Item(55, 12).Value = 600  # L55: 200 -> 600
$ws.Cells.Item(55, 14).Value = -946  # N55: -546 -> -946
$ws.Cells.Item(61, 8).Value = 4051.4375  # H61: 3043.75 -> 4051.4375
$ws.Cells.Item(61, 9).Value = 2568.5833  # I61: 2258 -> 2568.5833
$ws.Cells.Item(61, 10).Value = 8500  # J61: 5401 -> 8500
$ws.Cells.Item(61, 11).Value = 2568.5833  # K61: 2258 -> 2568.5833
$ws.Cells.Item(61, 12).Value = 8500  # L61: 5401 -> 8500
$ws.Cells.Item(61, 13).Value = -2366.5833  # M61: -2056 -> -2366.5833
$ws.Cells.Item(61, 14).Value = -8904  # N61: -5805 -> -8904
$ws.Cells.Item(68, 8).Value = 1583.3334  # H68: 1750 -> 1583.3334
$ws.Cells.Item(68, 10).Value = 1625  # J68: 2000 -> 1625
$ws.Cells.Item(68, 12).Value = 1625  # L68: 2000 -> 1625
$ws.Cells.Item(68, 14).Value = -3123  # N68: -3498 -> -3123
$ws.Cells.Item(71, 8).Value = 1583.3334  # H71: 1750 -> 1583.3334
$ws.Cells.Item(71, 10).Value = 1625  # J71: 2000 -> 1625
$ws.Cells.Item(71, 12).Value = 8125  # L71: 10000 -> 8125
$ws.Cells.Item(71, 14).Value = -15613  # N71: -17488 -> -15613
$ws.Cells.Item(113, 8).Value = 4051.4375  # H113: 3043.75 -> 4051.4375
$ws.Cells.Item(113, 9).Value = 2568.5833  # I113: 2258 -> 2568.5833
$ws.Cells.Item(113, 10).Value = 8500  # J113: 5401 -> 8500
$ws.Cells.Item(113, 11).Value = 2568.5833  # K113: 2258 -> 2568.5833
$ws.Cells.Item(113, 12).Value = 8500  # L113: 5401 -> 8500
$ws.Cells.Item(113, 13).Value = -398.5832999999998  # M113: -88 -> -398.5832999999998
$ws.Cells.Item(113, 14).Value = -12840  # N113: -9741 -> -12840
$ws.Cells.Item(126, 8).Value = 2005.8667  # H126: 2007.3077 -> 2005.8667
$ws.Cells.Item(126, 9).Value = 1807.6364  # I126: 2089.1 -> 1807.6364
$ws.Cells.Item(126, 10).Value = 2551  # J126: 1734.6666 -> 2551
$ws.Cells.Item(126, 11).Value = 5422.9092  # K126: 6267.299999999999 -> 5422.9092
$ws.Cells.Item(126, 12).Value = 7653  # L126: 5203.9998 -> 7653
$ws.Cells.Item(126, 13).Value = -2952.9092  # M126: -3797.299999999999 -> -2952.9092
$ws.Cells.Item(126, 14).Value = -12593  # N126: -10143.9998 -> -12593

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2432.4666  # H126: 2332.5 -> 2432.4666
$ws.Cells.Item(126, 9).Value = 2149.5  # I126: 2277.2222 -> 2149.5
$ws.Cells.Item(126, 10).Value = 2998.4  # J126: 2498.3333 -> 2998.4
$ws.Cells.Item(126, 11).Value = 6448.5  # K126: 6831.6666 -> 6448.5
$ws.Cells.Item(126, 12).Value = 8995.200000000001  # L126: 7494.999899999999 -> 8995.200000000001
$ws.Cells.Item(126, 13).Value = -3978.5  # M126: -4361.6666 -> -3978.5
$ws.Cells.Item(126, 14).Value = -13935.2  # N126: -12434.9999 -> -13935.2
